$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# The "Units" column values (unit1, unit2) for the two duplicate trait rows
# are being cleared out so the sample file no longer has duplicate trait
# names differentiated only by units.
$range = $ws.Range("M2:M3")
$range.Select() | Out-Null
$range.ClearContents() | Out-Null

$wb.Save()
